$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Aircraft is returning to the Home Point Minimum RTH Altitude is 30m You can reset the RTH Altitude in Remote Controller Settings after cancelling RTH ."
$ws.Range("C2").Value = "You can reset the RTH Altitude in Remote Controller Settings after cancelling RTH"
$ws.Range("D2").Value = "12-24"
$ws.Range("E2").Value = "NonEvent"
$ws.Range("F2").Value = "NonEvent"

# Row 3
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = "GPS signal weak Positioning accuracy affected Fly with caution ."
$ws.Range("C3").Value = "Fly with caution"
$ws.Range("D3").Value = "6-8"
$ws.Range("E3").Value = "NonEvent"
$ws.Range("F3").Value = "NonEvent"

# Row 4
$ws.Range("A4").Value = 11
$ws.Range("B4").Value = "Compass Error Compass data error Please check the Compass installation and interference ."
$ws.Range("C4").Value = "Please check the Compass installation and interference"
$ws.Range("D4").Value = "5-11"
$ws.Range("E4").Value = "NonEvent"
$ws.Range("F4").Value = "NonEvent"

# Row 5
$ws.Range("A5").Value = 12
$ws.Range("B5").Value = "Weak signal Avoid blocking the antennas and keep the antennas parallel to and facing toward the aircraft during flight Downlink Lost ."
$ws.Range("C5").Value = "Weak signal"
$ws.Range("D5").Value = "0-1"
$ws.Range("E5").Value = "Event"
$ws.Range("F5").Value = "Event"

# Row 6
$ws.Range("A6").Value = 15
$ws.Range("B6").Value = "High altitude Aircraft braking distance increased and flight time decreased Fly with caution ."
$ws.Range("C6").Value = "Fly with caution"
$ws.Range("D6").Value = "10-12"
$ws.Range("E6").Value = "NonEvent"
$ws.Range("F6").Value = "NonEvent"
